# ToDo.xlsx update: "Added mouse device to input manager. Implemented mouse
# controls in model viewer."
#
# - Removes the "Add input manager with keyboard and mouse" task row (old
#   row 2) entirely - the backlog item is done, so the whole row (and its
#   now-unused shared string) goes away and every row below shifts up one.
# - Revises the estimate on the "Add camera controls to model viewer" row
#   (old row 3 / new row 2) from 3 down to 2, now that mouse support has
#   partly landed.
# - Notes a new bug found while implementing it: the camera rolls on its
#   own, recorded as a comment on that same row's Task cell (new B2).
# - The three pre-existing comments live on cells that shifted up by one
#   row along with their data, so they're recreated at B4, B12 and B17
#   (previously B5, B13, B18) with their original text preserved.
# - Updates the active selection to the task cell that was annotated (B2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the text of the three existing comments before the row shift
# moves the cells they are anchored to.
$commentPointLights  = $ws.Range("B5").Comment.Text()
$commentShaders      = $ws.Range("B13").Comment.Text()
$commentModelCompiler = $ws.Range("B18").Comment.Text()

$ws.Range("B5").Comment.Delete()
$ws.Range("B13").Comment.Delete()
$ws.Range("B18").Comment.Delete()

# Delete the whole "Add input manager with keyboard and mouse" row; this is
# now superseded, everything below shifts up one row.
$ws.Rows(2).Delete()

# "Add camera controls to model viewer" (now row 2) - revised estimate.
$ws.Range("C2").Value = 2

# Recreate the three old comments at their new (shifted) cells.
$ws.Range("B4").AddComment($commentPointLights)
$ws.Range("B12").AddComment($commentShaders)
$ws.Range("B17").AddComment($commentModelCompiler)

# New comment documenting the camera bug found while adding mouse controls.
$ws.Range("B2").AddComment("Jonny:`nThe camera inadvertantly rolls around :-(")

# Match the author's final selection.
[void]$ws.Range("B2").Select()
